# Closed all risk items as they are no longer relevant to the project.
#
# The "Risk_Tracking_Log" sheet tracks individual risk rows in rows 7-15.
# Column B ("Current Status") drives an IF() formula in column E ("Risk
# Map") that reports "Closed" once the status is "Closed". Closing every
# risk means setting every row's status to "Closed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk_Tracking_Log")
$ws.Activate()

# Mirror the original author's editing session: start from the header
# row like before, then land the final selection on the last edited
# status cell (B15), matching the saved view state.
$ws.Range("A6").Select()

$ws.Range("B7:B15").Value = "Closed"

$ws.Range("B15").Select()

$wb.Save()
